# Finalizando importacion gestion y soporte
# - The "2: No" placeholder in the Si/No indicator legends (column C) is
#   corrected to "0: No" (matching the "0 or blank = NO" convention used by
#   the corresponding error-message cells) for each of the 9 affected rows.
# - The active cell / selection on the single worksheet is set to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("C19", "C22", "C32", "C37", "C42", "C43", "C44", "C45", "C52")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $text = $rng.Value2
    $rng.Value = $text.Replace("2: No", "0: No")
}

$ws.Range("B3").Select()
